# Restore C10 on the "Rules" sheet from 18 to 1 (per commit:
# "Restored from revision of admin on 10/27/2020 08:09:25 AM.TEST").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
